$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.225.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.878.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4309"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3708"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8846"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.908.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.500"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06990"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009153"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.28%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.266.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("E23").Value = "  +3.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.124.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.977"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.440"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.878"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08981"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7967"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.729"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.190"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.948"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.132"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.000"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05468"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01969"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.888"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5185"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1696"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.898"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.630"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4779"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06588"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.659"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.863"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.45%  "
